$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (preserve trailing zeros / exact formatting)
$textCells = @('D4', 'D5', 'D8', 'D11', 'D14', 'D15', 'D17', 'D18', 'D21', 'D23', 'D24', 'D25', 'D28', 'D29', 'D31', 'D33', 'D39', 'D40', 'D43', 'D44', 'D45', 'D46', 'D48', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '27.904.71'
$ws.Range('E2').Value = '  +2.98%  '
$ws.Range('D3').Value = '1.571.27'
$ws.Range('E3').Value = '  +0.40%  '
$ws.Range('D4').Value = '0.990'
$ws.Range('E4').Value = '  -1.78%  '
$ws.Range('D5').Value = '211.48'
$ws.Range('E5').Value = '  +0.53%  '
$ws.Range('E6').Value = '  +0.98%  '
$ws.Range('E7').Value = '  -1.83%  '
$ws.Range('D8').Value = '23.29'
$ws.Range('E8').Value = '  +6.32%  '
$ws.Range('E9').Value = '  +0.74%  '
$ws.Range('E10').Value = '  -0.02%  '
$ws.Range('D11').Value = '0.0877'
$ws.Range('E11').Value = '  +1.83%  '
$ws.Range('D12').Value = '1.794.24'
$ws.Range('E12').Value = '  +0.33%  '
$ws.Range('D13').Value = '1.571.56'
$ws.Range('E13').Value = '  +0.42%  '
$ws.Range('D14').Value = '3.76'
$ws.Range('E14').Value = '  -0.39%  '
$ws.Range('D15').Value = '0.521'
$ws.Range('E15').Value = '  +0.46%  '
$ws.Range('D16').Value = '27.870.92'
$ws.Range('E16').Value = '  +2.87%  '
$ws.Range('D17').Value = '63.46'
$ws.Range('E17').Value = '  +2.48%  '
$ws.Range('D18').Value = '230.59'
$ws.Range('E18').Value = '  +7.42%  '
$ws.Range('D19').Value = '0.0₃0704'
$ws.Range('E19').Value = '  +0.27%  '
$ws.Range('E20').Value = '  +1.22%  '
$ws.Range('D21').Value = '0.990'
$ws.Range('E21').Value = '  -1.80%  '
$ws.Range('E22').Value = '  -0.47%  '
$ws.Range('D23').Value = '9.31'
$ws.Range('E23').Value = '  +1.26%  '
$ws.Range('D24').Value = '1.93'
$ws.Range('E24').Value = '  -0.23%  '
$ws.Range('D25').Value = '151.19'
$ws.Range('E25').Value = '  -1.88%  '
$ws.Range('E26').Value = '  +1.39%  '
$ws.Range('E27').Value = '  +1.28%  '
$ws.Range('D28').Value = '6.57'
$ws.Range('E28').Value = '  -0.42%  '
$ws.Range('D29').Value = '0.991'
$ws.Range('E29').Value = '  -1.63%  '
$ws.Range('E30').Value = '  -0.11%  '
$ws.Range('D31').Value = '0.0473'
$ws.Range('E31').Value = '  +0.24%  '
$ws.Range('E32').Value = '  -0.22%  '
$ws.Range('D33').Value = '3.14'
$ws.Range('E33').Value = '  -1.51%  '
$ws.Range('D34').Value = '1.418.55'
$ws.Range('E34').Value = '  -0.56%  '
$ws.Range('E35').Value = '  -1.15%  '
$ws.Range('E36').Value = '  -4.92%  '
$ws.Range('E37').Value = '  -2.28%  '
$ws.Range('E38').Value = '  +0.14%  '
$ws.Range('D39').Value = '0.543'
$ws.Range('E39').Value = '  +2.64%  '
$ws.Range('D40').Value = '2.41'
$ws.Range('E40').Value = '  +2.00%  '
$ws.Range('E41').Value = '  -0.07%  '
$ws.Range('E42').Value = '  -1.81%  '
$ws.Range('D43').Value = '5.59'
$ws.Range('E43').Value = '  -3.93%  '
$ws.Range('D44').Value = '1.81'
$ws.Range('E44').Value = '  +4.41%  '
$ws.Range('D45').Value = '0.966'
$ws.Range('E45').Value = '  -3.59%  '
$ws.Range('D46').Value = '63.97'
$ws.Range('E46').Value = '  -0.79%  '
$ws.Range('D47').Value = '1.704.79'
$ws.Range('E47').Value = '  -0.29%  '
$ws.Range('D48').Value = '86.86'
$ws.Range('E48').Value = '  +1.04%  '
$ws.Range('E49').Value = '  +1.11%  '
$ws.Range('D50').Value = '0.0₆0101'
$ws.Range('E50').Value = '  -1.11%  '
$ws.Range('D51').Value = '39.49'
$ws.Range('E51').Value = '  +16.67%  '
